$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.032022510714174
$ws.Range("D2").Value = 1.039439945612714
$ws.Range("E2").Value = 1.035662548159038
$ws.Range("F2").Value = 1.04665282801693
$ws.Range("I2").Value = 1.036623732717049
$ws.Range("J2").Value = 1.037155017940469
$ws.Range("K2").Value = 1.042225128265613
$ws.Range("L2").Value = 1.038458509889307
$ws.Range("M2").Value = 1.04941766266144
$ws.Range("N2").Value = 1.016362412607031

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.03283848813016
$ws.Range("D3").Value = 1.040072315802218
$ws.Range("E3").Value = 1.036426211215169
$ws.Range("F3").Value = 1.047466279922317
$ws.Range("I3").Value = 1.036795366223756
$ws.Range("J3").Value = 1.037613794783101
$ws.Range("K3").Value = 1.042668186436973
$ws.Range("L3").Value = 1.039031735889701
$ws.Range("M3").Value = 1.050042795848377
$ws.Range("N3").Value = 1.016515960411325

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.033367174661673
$ws.Range("D4").Value = 1.040482116614585
$ws.Range("E4").Value = 1.036921382031249
$ws.Range("F4").Value = 1.047993732321753
$ws.Range("I4").Value = 1.036905570304299
$ws.Range("J4").Value = 1.037910680313308
$ws.Range("K4").Value = 1.042954799093532
$ws.Range("L4").Value = 1.03940301887937
$ws.Range("M4").Value = 1.050447739748749
$ws.Range("N4").Value = 1.016615284520294

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.033589599032813
$ws.Range("D5").Value = 1.0406545422796
$ws.Range("E5").Value = 1.037129796432957
$ws.Range("F5").Value = 1.04821573333743
$ws.Range("I5").Value = 1.036951694866977
$ws.Range("J5").Value = 1.038035495979661
$ws.Range("K5").Value = 1.043075271620341
$ws.Range("L5").Value = 1.039559192667408
$ws.Range("M5").Value = 1.050618082039646
$ws.Range("N5").Value = 1.01665703233633

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.033626954662672
$ws.Range("D6").Value = 1.040683501780687
$ws.Range("E6").Value = 1.037164804414722
$ws.Range("F6").Value = 1.048253023473954
$ws.Range("I6").Value = 1.036959427328232
$ws.Range("J6").Value = 1.0380564533343
$ws.Range("K6").Value = 1.043095498309305
$ws.Range("L6").Value = 1.03958541995633
$ws.Range("M6").Value = 1.05064668930217
$ws.Range("N6").Value = 1.016664041493028

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.033370146062589
$ws.Range("D7").Value = 1.040484420005729
$ws.Range("E7").Value = 1.036924165916043
$ws.Range("F7").Value = 1.04799669769101
$ws.Range("I7").Value = 1.036906187430281
$ws.Range("J7").Value = 1.037912348088366
$ws.Range("K7").Value = 1.042956408931368
$ws.Range("L7").Value = 1.039405105343542
$ws.Range("M7").Value = 1.050450015464733
$ws.Range("N7").Value = 1.016615842388845

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.032298129243325
$ws.Range("D8").Value = 1.039653529635137
$ws.Range("E8").Value = 1.035920417283924
$ws.Range("F8").Value = 1.046927510028791
$ws.Range("I8").Value = 1.036681913515568
$ws.Range("J8").Value = 1.037310057513658
$ws.Range("K8").Value = 1.042374876702046
$ws.Range("L8").Value = 1.038652157290242
$ws.Range("M8").Value = 1.049628837039676
$ws.Range("N8").Value = 1.016414311060057

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.03041449495762
$ws.Range("D9").Value = 1.038194185923415
$ws.Range("E9").Value = 1.034159656215723
$ws.Range("F9").Value = 1.04505193799485
$ws.Range("I9").Value = 1.036280203305077
$ws.Range("J9").Value = 1.036249002639929
$ws.Range("K9").Value = 1.041349622558983
$ws.Range("L9").Value = 1.037328249141949
$ws.Range("M9").Value = 1.048185267205029
$ws.Range("N9").Value = 1.016058966282134

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.029162470320995
$ws.Range("D10").Value = 1.037224621006888
$ws.Range("E10").Value = 1.032991291255478
$ws.Range("F10").Value = 1.043807374826085
$ws.Range("I10").Value = 1.036008061906288
$ws.Range("J10").Value = 1.03554188339682
$ws.Range("K10").Value = 1.040665849853208
$ws.Range("L10").Value = 1.036447674146314
$ws.Range("M10").Value = 1.047225308542532
$ws.Range("N10").Value = 1.015821949782445

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.028621236711374
$ws.Range("D11").Value = 1.036805601934494
$ws.Range("E11").Value = 1.032486699591943
$ws.Range("F11").Value = 1.043269870907488
$ws.Range("I11").Value = 1.035889202649187
$ws.Range("J11").Value = 1.03523576872494
$ws.Range("K11").Value = 1.040369721556383
$ws.Range("L11").Value = 1.036066876589112
$ws.Range("M11").Value = 1.046810231427294
$ws.Range("N11").Value = 1.015719296586966

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.028420335620949
$ws.Range("D12").Value = 1.036650083013782
$ws.Range("E12").Value = 1.032299471799216
$ws.Range("F12").Value = 1.043070430464626
$ws.Range("I12").Value = 1.035844900470971
$ws.Range("J12").Value = 1.035122076333331
$ws.Range("K12").Value = 1.040259720140187
$ws.Range("L12").Value = 1.035925507706749
$ws.Range("M12").Value = 1.04665614387067
$ws.Range("N12").Value = 1.015681163576975

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.028463423351311
$ws.Range("D13").Value = 1.036683436734939
$ws.Range("E13").Value = 1.032339623713101
$ws.Range("F13").Value = 1.043113201469529
$ws.Range("I13").Value = 1.03585441033056
$ws.Range("J13").Value = 1.035146463162829
$ws.Range("K13").Value = 1.040283316075256
$ws.Range("L13").Value = 1.035955828334422
$ws.Range("M13").Value = 1.046689192048648
$ws.Range("N13").Value = 1.015689343366945

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.028604627356153
$ws.Range("D14").Value = 1.036792744163371
$ws.Range("E14").Value = 1.032471219192364
$ws.Range("F14").Value = 1.043253380752377
$ws.Range("I14").Value = 1.035885543725274
$ws.Range("J14").Value = 1.035226370619339
$ws.Range("K14").Value = 1.040360628923238
$ws.Range("L14").Value = 1.036055189427844
$ws.Range("M14").Value = 1.04679749264431
$ws.Range("N14").Value = 1.015716144559721

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.028691646032922
$ws.Range("D15").Value = 1.036860108488305
$ws.Range("E15").Value = 1.032552326074246
$ws.Range("F15").Value = 1.043339778048513
$ws.Range("I15").Value = 1.035904705843484
$ws.Range("J15").Value = 1.035275605900174
$ws.Range("K15").Value = 1.040408263142103
$ws.Range("L15").Value = 1.036116419232724
$ws.Range("M15").Value = 1.046864232268074
$ws.Range("N15").Value = 1.015732657270508

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.029198409352875
$ws.Range("D16").Value = 1.037252447127109
$ws.Range("E16").Value = 1.033024807282749
$ws.Range("F16").Value = 1.043843076868714
$ws.Range("I16").Value = 1.036015928770033
$ws.Range("J16").Value = 1.035562200861912
$ws.Range("K16").Value = 1.040685501971432
$ws.Range("L16").Value = 1.036472957037159
$ws.Range("M16").Value = 1.047252868449429
$ws.Range("N16").Value = 1.015828762083605

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.029516531143558
$ws.Range("D17").Value = 1.037498768696882
$ws.Range("E17").Value = 1.033321536554695
$ws.Range("F17").Value = 1.044159159255211
$ws.Range("I17").Value = 1.036085423335687
$ws.Range("J17").Value = 1.035741994682778
$ws.Range("K17").Value = 1.040859394084377
$ws.Range("L17").Value = 1.036696738021345
$ws.Range("M17").Value = 1.047496809190147
$ws.Range("N17").Value = 1.015889040148927

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.029702173048155
$ws.Range("D18").Value = 1.037642521825181
$ws.Range("E18").Value = 1.033494740760259
$ws.Range("F18").Value = 1.044343659717333
$ws.Range("I18").Value = 1.036125859853688
$ws.Range("J18").Value = 1.035846872295773
$ws.Range("K18").Value = 1.040960817476714
$ws.Range("L18").Value = 1.036827313552396
$ws.Range("M18").Value = 1.047639152635148
$ws.Range("N18").Value = 1.015924197025866

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.029765486817425
$ws.Range("D19").Value = 1.037691551086873
$ws.Range("E19").Value = 1.033553820406176
$ws.Range("F19").Value = 1.044406592433948
$ws.Range("I19").Value = 1.036139630934448
$ws.Range("J19").Value = 1.035882633985101
$ws.Range("K19").Value = 1.040995399354005
$ws.Range("L19").Value = 1.036871844505379
$ws.Range("M19").Value = 1.047687697676585
$ws.Range("N19").Value = 1.015936184207667

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.0294823906905
$ws.Range("D20").Value = 1.037472332649951
$ws.Range("E20").Value = 1.033289687157167
$ws.Range("F20").Value = 1.044125232632551
$ws.Range("I20").Value = 1.036077977408197
$ws.Range("J20").Value = 1.03572270378218
$ws.Range("K20").Value = 1.040840737607006
$ws.Range("L20").Value = 1.03667272347685
$ws.Range("M20").Value = 1.047470630755479
$ws.Range("N20").Value = 1.015882573117247

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.028563042495433
$ws.Range("D21").Value = 1.036760552427991
$ws.Range("E21").Value = 1.03243246206618
$ws.Range("F21").Value = 1.043212095566165
$ws.Range("I21").Value = 1.035876379922056
$ws.Range("J21").Value = 1.035202839516935
$ws.Range("K21").Value = 1.040337862376195
$ws.Range("L21").Value = 1.036025927943782
$ws.Range("M21").Value = 1.046765598304302
$ws.Range("N21").Value = 1.015708252356538

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.027985806295634
$ws.Range("D22").Value = 1.036313742951722
$ws.Range("E22").Value = 1.031894648856522
$ws.Range("F22").Value = 1.042639200224234
$ws.Range("I22").Value = 1.035748745286295
$ws.Range("J22").Value = 1.034876051614883
$ws.Range("K22").Value = 1.040021649668001
$ws.Range("L22").Value = 1.035619704067923
$ws.Range("M22").Value = 1.046322841256686
$ws.Range("N22").Value = 1.015598632612571

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.028291734229743
$ws.Range("D23").Value = 1.036550536706041
$ws.Range("E23").Value = 1.032179643363083
$ws.Range("F23").Value = 1.042942785580932
$ws.Range("I23").Value = 1.035816490211766
$ws.Range("J23").Value = 1.035049280857906
$ws.Range("K23").Value = 1.040189282879452
$ws.Range("L23").Value = 1.035835008591675
$ws.Range("M23").Value = 1.04655750482005
$ws.Range("N23").Value = 1.015656745596506

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.029497817023558
$ws.Range("D24").Value = 1.037484277721111
$ws.Range("E24").Value = 1.033304078136695
$ws.Range("F24").Value = 1.044140562196268
$ws.Range("I24").Value = 1.036081342206722
$ws.Range("J24").Value = 1.035731420488583
$ws.Range("K24").Value = 1.040849167681496
$ws.Range("L24").Value = 1.036683574468112
$ws.Range("M24").Value = 1.047482459487564
$ws.Range("N24").Value = 1.015885495297659

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.030900809440055
$ws.Range("D25").Value = 1.03857088217098
$ws.Range("E25").Value = 1.03461389848321
$ws.Range("F25").Value = 1.04553580170532
$ws.Range("I25").Value = 1.036384822179073
$ws.Range("J25").Value = 1.036523272226108
$ws.Range("K25").Value = 1.041614728746268
$ws.Range("L25").Value = 1.037670160018253
$ws.Range("M25").Value = 1.048558045221443
$ws.Range("N25").Value = 1.016150854657187
